$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying source data had rows reshuffled within each match-day block
# (same Div/Date), so we rotate/swap the B:B and E:AD content among the rows
# of each block while keeping column A (sequential id) fixed in place.

# --- Block: rows 15, 16, 17 ---
$b15 = $ws.Range("B15").Value()
$ead15 = $ws.Range("E15:AD15").Value()
$b16 = $ws.Range("B16").Value()
$ead16 = $ws.Range("E16:AD16").Value()
$b17 = $ws.Range("B17").Value()
$ead17 = $ws.Range("E17:AD17").Value()

$ws.Range("B15").Value = $b17
$ws.Range("E15:AD15").Value = $ead17
$ws.Range("B16").Value = $b15
$ws.Range("E16:AD16").Value = $ead15
$ws.Range("B17").Value = $b16
$ws.Range("E17:AD17").Value = $ead16

# --- Block: rows 131, 132 ---
$b131 = $ws.Range("B131").Value()
$ead131 = $ws.Range("E131:AD131").Value()
$b132 = $ws.Range("B132").Value()
$ead132 = $ws.Range("E132:AD132").Value()

$ws.Range("B131").Value = $b132
$ws.Range("E131:AD131").Value = $ead132
$ws.Range("B132").Value = $b131
$ws.Range("E132:AD132").Value = $ead131

# --- Block: rows 136, 137 ---
$b136 = $ws.Range("B136").Value()
$ead136 = $ws.Range("E136:AD136").Value()
$b137 = $ws.Range("B137").Value()
$ead137 = $ws.Range("E137:AD137").Value()

$ws.Range("B136").Value = $b137
$ws.Range("E136:AD136").Value = $ead137
$ws.Range("B137").Value = $b136
$ws.Range("E137:AD137").Value = $ead136

# --- Block: rows 144, 145 ---
$b144 = $ws.Range("B144").Value()
$ead144 = $ws.Range("E144:AD144").Value()
$b145 = $ws.Range("B145").Value()
$ead145 = $ws.Range("E145:AD145").Value()

$ws.Range("B144").Value = $b145
$ws.Range("E144:AD144").Value = $ead145
$ws.Range("B145").Value = $b144
$ws.Range("E145:AD145").Value = $ead144

# --- Block: rows 159, 160 ---
$b159 = $ws.Range("B159").Value()
$ead159 = $ws.Range("E159:AD159").Value()
$b160 = $ws.Range("B160").Value()
$ead160 = $ws.Range("E160:AD160").Value()

$ws.Range("B159").Value = $b160
$ws.Range("E159:AD159").Value = $ead160
$ws.Range("B160").Value = $b159
$ws.Range("E160:AD160").Value = $ead159

# --- Block: rows 176, 177, 178, 179, 180 ---
$b176 = $ws.Range("B176").Value()
$ead176 = $ws.Range("E176:AD176").Value()
$b177 = $ws.Range("B177").Value()
$ead177 = $ws.Range("E177:AD177").Value()
$b178 = $ws.Range("B178").Value()
$ead178 = $ws.Range("E178:AD178").Value()
$b179 = $ws.Range("B179").Value()
$ead179 = $ws.Range("E179:AD179").Value()
$b180 = $ws.Range("B180").Value()
$ead180 = $ws.Range("E180:AD180").Value()

$ws.Range("B176").Value = $b178
$ws.Range("E176:AD176").Value = $ead178
$ws.Range("B177").Value = $b180
$ws.Range("E177:AD177").Value = $ead180
$ws.Range("B178").Value = $b177
$ws.Range("E178:AD178").Value = $ead177
$ws.Range("B179").Value = $b176
$ws.Range("E179:AD179").Value = $ead176
$ws.Range("B180").Value = $b179
$ws.Range("E180:AD180").Value = $ead179

# --- Block: rows 200, 202 ---
$b200 = $ws.Range("B200").Value()
$ead200 = $ws.Range("E200:AD200").Value()
$b202 = $ws.Range("B202").Value()
$ead202 = $ws.Range("E202:AD202").Value()

$ws.Range("B200").Value = $b202
$ws.Range("E200:AD200").Value = $ead202
$ws.Range("B202").Value = $b200
$ws.Range("E202:AD202").Value = $ead200
